$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-24 05:18:11"
$ws.Range("N2").Value = "1.0 °C 4:59 TU"
$ws.Range("E3").Value = "2026-02-24 05:18:13"
$ws.Range("K3").Value = "-0.1 MJ/m2"
$ws.Range("N3").Value = "0.5 °C 4:50 TU"
$ws.Range("O3").Value = "2.7 °C"
$ws.Range("E4").Value = "2026-02-24 05:18:16"
$ws.Range("J4").Value = "1022.5 hPa"
$ws.Range("N4").Value = "5.1 °C 4:50 TU"
$ws.Range("O4").Value = "6.8 °C"
$ws.Range("E5").Value = "2026-02-24 05:18:18"
$ws.Range("K5").Value = "-0.1 MJ/m2"
$ws.Range("E6").Value = "2026-02-24 05:18:20"
$ws.Range("J6").Value = "1022.0 hPa"
$ws.Range("N6").Value = "8.0 °C 4:43 TU"
$ws.Range("O6").Value = "9.5 °C"
$ws.Range("E7").Value = "2026-02-24 05:18:23"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "82%"
$ws.Range("C7").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("J7").Value = "1022.2 hPa"
$ws.Range("N7").Value = "11.3 °C 4:59 TU"
$ws.Range("O7").Value = "12.2 °C"
$ws.Range("E8").Value = "2026-02-24 05:18:25"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "45%"
$ws.Range("C8").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("J8").Value = "1021.8 hPa"
$ws.Range("E9").Value = "2026-02-24 05:18:27"
$ws.Range("O9").Value = "5.6 °C"
$ws.Range("E10").Value = "2026-02-24 05:18:30"
$ws.Range("K10").Value = "-0.1 MJ/m2"
$ws.Range("E11").Value = "2026-02-24 05:18:32"
$ws.Range("N11").Value = "1.5 °C 4:59 TU"
$ws.Range("O11").Value = "2.7 °C"
$ws.Range("E12").Value = "2026-02-24 05:18:34"
$ws.Range("N12").Value = "3.6 °C 4:50 TU"
$ws.Range("O12").Value = "6.1 °C"
$ws.Range("E13").Value = "2026-02-24 05:18:36"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "92%"
$ws.Range("C13").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("K13").Value = "-0.1 MJ/m2"
$ws.Range("N13").Value = "-3.3 °C 4:32 TU"
$ws.Range("O13").Value = "-1.3 °C"
$ws.Range("E14").Value = "2026-02-24 05:18:39"
$ws.Range("N14").Value = "8.0 °C 4:36 TU"
$ws.Range("E15").Value = "2026-02-24 05:18:41"
$ws.Range("N15").Value = "4.6 °C 4:59 TU"
$ws.Range("O15").Value = "6.0 °C"
$ws.Range("E16").Value = "2026-02-24 05:18:43"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "20%"
$ws.Range("C16").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("N16").Value = "3.1 °C 4:47 TU"
$ws.Range("E17").Value = "2026-02-24 05:18:46"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "30%"
$ws.Range("C17").Copy()
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("K17").Value = "-0.1 MJ/m2"
$ws.Range("E18").Value = "2026-02-24 05:18:48"
$ws.Range("J18").Value = "1022.7 hPa"
$ws.Range("N18").Value = "2.0 °C 4:59 TU"
$ws.Range("O18").Value = "3.1 °C"
$ws.Range("E19").Value = "2026-02-24 05:18:51"
$ws.Range("K19").Value = "-0.1 MJ/m2"
$ws.Range("E20").Value = "2026-02-24 05:18:53"
$ws.Range("N20").Value = "0.1 °C 4:39 TU"
$ws.Range("O20").Value = "1.0 °C"
$ws.Range("E21").Value = "2026-02-24 05:18:55"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "80%"
$ws.Range("C21").Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("K21").Value = "-0.1 MJ/m2"
$ws.Range("N21").Value = "2.0 °C 4:59 TU"
$ws.Range("O21").Value = "3.7 °C"
$ws.Range("E22").Value = "2026-02-24 05:18:58"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "18%"
$ws.Range("C22").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("E23").Value = "2026-02-24 05:19:00"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "28%"
$ws.Range("C23").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("K23").Value = "-0.1 MJ/m2"
$ws.Range("E24").Value = "2026-02-24 05:19:02"
$ws.Range("N24").Value = "1.1 °C 4:57 TU"
$ws.Range("O24").Value = "3.4 °C"
$ws.Range("E25").Value = "2026-02-24 05:19:05"
$ws.Range("K25").Value = "-0.1 MJ/m2"
$ws.Range("E26").Value = "2026-02-24 05:19:07"
$ws.Range("J26").Value = "1022.7 hPa"
$ws.Range("K26").Value = "-0.1 MJ/m2"
$ws.Range("L26").Value = "11.9 km/h - 13º 4:43 TU"
$ws.Range("N26").Value = "6.5 °C 4:45 TU"
$ws.Range("E27").Value = "2026-02-24 05:19:10"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "39%"
$ws.Range("C27").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("K27").Value = "-0.1 MJ/m2"
$ws.Range("E28").Value = "2026-02-24 05:19:12"
$ws.Range("J28").Value = "1024.0 hPa"
$ws.Range("O28").Value = "3.8 °C"
$ws.Range("E29").Value = "2026-02-24 05:19:14"
$ws.Range("N29").Value = "3.5 °C 4:59 TU"
$ws.Range("O29").Value = "5.2 °C"
$ws.Range("E30").Value = "2026-02-24 05:19:17"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "82%"
$ws.Range("C30").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("J30").Value = "1022.0 hPa"
$ws.Range("N30").Value = "7.9 °C 4:57 TU"
$ws.Range("O30").Value = "9.2 °C"
$ws.Range("E31").Value = "2026-02-24 05:19:19"
$ws.Range("J31").Value = "1021.1 hPa"
$ws.Range("N31").Value = "13.8 °C 4:59 TU"
$ws.Range("E32").Value = "2026-02-24 05:19:21"
$ws.Range("K32").Value = "-0.1 MJ/m2"
$ws.Range("N32").Value = "-4.1 °C 4:30 TU"
$ws.Range("O32").Value = "-3.1 °C"
$ws.Range("E33").Value = "2026-02-24 05:19:24"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "69%"
$ws.Range("C33").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("N33").Value = "0.4 °C 4:51 TU"
$ws.Range("O33").Value = "2.2 °C"
$ws.Range("E34").Value = "2026-02-24 05:19:26"
$ws.Range("O34").Value = "2.7 °C"
$ws.Range("E35").Value = "2026-02-24 05:19:29"
$ws.Range("E36").Value = "2026-02-24 05:19:31"
$ws.Range("J36").Value = "1021.9 hPa"
$ws.Range("E37").Value = "2026-02-24 05:19:33"
$ws.Range("J37").Value = "1027.6 hPa"
$ws.Range("N37").Value = "-0.3 °C 4:59 TU"
$ws.Range("O37").Value = "0.8 °C"
$ws.Range("E38").Value = "2026-02-24 05:19:36"
$ws.Range("K38").Value = "-0.1 MJ/m2"
$ws.Range("O38").Value = "6.8 °C"
$ws.Range("E39").Value = "2026-02-24 05:19:38"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "34%"
$ws.Range("C39").Copy()
$ws.Range("H39").PasteSpecial(-4122)
$ws.Range("E40").Value = "2026-02-24 05:19:40"
$ws.Range("N40").Value = "0.1 °C 4:59 TU"
$ws.Range("O40").Value = "1.3 °C"
$ws.Range("E41").Value = "2026-02-24 05:19:43"
$ws.Range("J41").Value = "1022.3 hPa"
$ws.Range("E42").Value = "2026-02-24 05:19:45"
$ws.Range("O42").Value = "6.9 °C"
$ws.Range("E43").Value = "2026-02-24 05:19:47"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "99%"
$ws.Range("C43").Copy()
$ws.Range("H43").PasteSpecial(-4122)
$ws.Range("N43").Value = "2.5 °C 4:39 TU"
$ws.Range("O43").Value = "4.2 °C"
$ws.Range("E44").Value = "2026-02-24 05:19:50"
$ws.Range("K44").Value = "-0.1 MJ/m2"
$ws.Range("L44").Value = "6.8 km/h - 21º 4:59 TU"
$ws.Range("N44").Value = "-2.1 °C 4:41 TU"
$ws.Range("O44").Value = "0.0 °C"
$ws.Range("E45").Value = "2026-02-24 05:19:52"
$ws.Range("K45").Value = "-0.1 MJ/m2"
$ws.Range("L45").Value = "15.1 km/h - 85º 4:55 TU"
$ws.Range("O45").Value = "4.4 °C"
$ws.Range("E46").Value = "2026-02-24 05:19:55"
$ws.Range("K46").Value = "-0.1 MJ/m2"
$ws.Range("N46").Value = "1.0 °C 4:58 TU"
$ws.Range("O46").Value = "2.4 °C"
$excel.CutCopyMode = 0
